# Re-process the data with the newly curated dimensions: the
# "edad-grupos-quinquenales-2010" and "sexo" columns move from being
# dimensions to being measures, and their now-unused mapping-file rows
# (row 5) are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: edad-grupos-quinquenales-2010 dimension -> measure
$ws.Range("A2").Value = "iaest-measure:edad-grupos-quinquenales-2010"
$ws.Range("A3").Value = "medida"
$ws.Range("A4").Value = "xsd:int"

# Column F: sexo dimension -> measure
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"

# Row 5 held the mapping-file references for the two former dimensions;
# those no longer apply now that both columns are measures.
$ws.Rows.Item(5).Delete()
